$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TDOC")

# Row 16 - Gross Margin
$ws.Range("D16").Value = 0.5977
$ws.Range("E16").Value = 0.6165
$ws.Range("F16").Value = 0.6402
$ws.Range("G16").Value = 0.6666

# Row 20 - Free Cash Flow Margin
$ws.Range("D20").Value = 0.0728
$ws.Range("E20").Value = 0.0682
$ws.Range("F20").Value = 0.034
$ws.Range("G20").Value = 0.0343

# Row 28 - EBITDA Margin
$ws.Range("D28").Value = 0.0393
$ws.Range("E28").Value = 0.0113
$ws.Range("F28").Value = -0.0119
$ws.Range("G28").Value = -0.0182

# Row 29 - Operating Cash Flow Margin
$ws.Range("D29").Value = 0.0968
$ws.Range("E29").Value = 0.0888
$ws.Range("F29").Value = 0.054
$ws.Range("G29").Value = 0.054
